$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Note: Price values in column D (e.g. "67.774.03", "0.100") are stored as
# text, not numbers (the sheet uses "." as a thousands-style separator, and
# some values have significant trailing zeros). A leading apostrophe forces
# Excel to keep these assignments as literal text instead of auto-converting
# them to numeric/date values.
$ws.Range('D2').Value = '''67.774.03'
$ws.Range('E2').Value = '  -1.34%  '
$ws.Range('D3').Value = '''3.788.96'
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''593.87'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('D6').Value = '''166.48'
$ws.Range('E6').Value = '  -1.45%  '
$ws.Range('D7').Value = '''3.787.14'
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -0.61%  '
$ws.Range('E10').Value = '  -1.23%  '
$ws.Range('D11').Value = '''6.40'
$ws.Range('E11').Value = '  -1.70%  '
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('D13').Value = '''0.0000257'
$ws.Range('E13').Value = '  -3.37%  '
$ws.Range('E14').Value = '  -1.91%  '
$ws.Range('D15').Value = '''4.415.60'
$ws.Range('E15').Value = '  -0.07%  '
$ws.Range('D16').Value = '''3.772.68'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('D17').Value = '''67.863.94'
$ws.Range('E17').Value = '  -1.13%  '
$ws.Range('D18').Value = '''17.82'
$ws.Range('E18').Value = '  -2.84%  '
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('E20').Value = '  -2.01%  '
$ws.Range('D21').Value = '''10.56'
$ws.Range('E21').Value = '  -1.96%  '
$ws.Range('D22').Value = '''461.01'
$ws.Range('E22').Value = '  -1.90%  '
$ws.Range('D23').Value = '''0.697'
$ws.Range('E23').Value = '  -0.97%  '
$ws.Range('E24').Value = '  +5.55%  '
$ws.Range('D25').Value = '''83.64'
$ws.Range('E25').Value = '  -1.21%  '
$ws.Range('D26').Value = '''2.15'
$ws.Range('E26').Value = '  -5.00%  '
$ws.Range('D27').Value = '''11.86'
$ws.Range('E27').Value = '  -2.85%  '
$ws.Range('E28').Value = '  -2.64%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  -2.24%  '
$ws.Range('D31').Value = '''29.87'
$ws.Range('E31').Value = '  -1.14%  '
$ws.Range('D32').Value = '''7.21'
$ws.Range('E32').Value = '  -3.58%  '
$ws.Range('E33').Value = '  -3.59%  '
$ws.Range('E34').Value = '  +0.32%  '
$ws.Range('D35').Value = '''9.10'
$ws.Range('E35').Value = '  -1.61%  '
$ws.Range('D36').Value = '''3.736.14'
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('D37').Value = '''0.100'
$ws.Range('E37').Value = '  -2.19%  '
$ws.Range('D38').Value = '''3.46'
$ws.Range('E38').Value = '  -0.94%  '
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('D40').Value = '''0.996'
$ws.Range('E40').Value = '  -0.72%  '
$ws.Range('E41').Value = '  -1.64%  '
$ws.Range('D42').Value = '''0.999'
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('D44').Value = '''43.77'
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').Value = '''0.299'
$ws.Range('E45').Value = '  -3.14%  '
$ws.Range('D46').Value = '''46.77'
$ws.Range('E46').Value = '  +2.77%  '
$ws.Range('E47').Value = '  -4.61%  '
$ws.Range('E48').Value = '  -3.06%  '
$ws.Range('D49').Value = '''146.29'
$ws.Range('E49').Value = '  +0.69%  '
$ws.Range('D50').Value = '''386.55'
$ws.Range('E50').Value = '  -5.34%  '
$ws.Range('D51').Value = '''2.758.34'
$ws.Range('E51').Value = '  +3.48%  '
